$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: volume number and report week dates ---
# "Volume 30   Number  29" -> "Volume 30   Number  30"
$ws.Range("A8").Characters(21, 2).Text = "30"
# "Report Covering the Week  7/17/2023  Through  7/23/2023"
#   -> "Report Covering the Week  7/24/2023  Through  7/30/2023"
$ws.Range("C9").Characters(27, 9).Text = "7/24/2023"
$ws.Range("C9").Characters(47, 9).Text = "7/30/2023"

# --- Crime statistics table (rows 15-30) ---

# Row 15 (Rape)
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = "'0"
$ws.Range("A15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = "'***.*"
$ws.Range("A15").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("I15").Value = 13
$ws.Range("K15").Value = 0
$ws.Range("M15").Value = 85.714285714285
$ws.Range("N15").Value = -50

# Row 16 (Robbery)
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = -19.047619047619
$ws.Range("I16").Value = 106
$ws.Range("J16").Value = 152
$ws.Range("K16").Value = -30.263157894736
$ws.Range("L16").Value = 30.864197530864
$ws.Range("M16").Value = 4.950495049504
$ws.Range("N16").Value = -79.497098646034

# Row 17 (Fel. Assault)
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = -31.578947368421
$ws.Range("I17").Value = 141
$ws.Range("J17").Value = 111
$ws.Range("K17").Value = 27.027027027027
$ws.Range("L17").Value = 36.893203883495
$ws.Range("M17").Value = 56.666666666666
$ws.Range("N17").Value = -59.482758620689

# Row 18 (Burglary)
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -28.571428571428
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 32
$ws.Range("H18").Value = -53.125
$ws.Range("I18").Value = 157
$ws.Range("J18").Value = 227
$ws.Range("K18").Value = -30.837004405286
$ws.Range("L18").Value = -10.285714285714
$ws.Range("M18").Value = 12.949640287769
$ws.Range("N18").Value = -69.033530571992

# Row 19 (Gr. Larceny)
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 33
$ws.Range("E19").Value = -75.757575757575
$ws.Range("F19").Value = 76
$ws.Range("G19").Value = 96
$ws.Range("H19").Value = -20.833333333333
$ws.Range("I19").Value = 574
$ws.Range("J19").Value = 577
$ws.Range("K19").Value = -0.519930675909
$ws.Range("L19").Value = 57.260273972602
$ws.Range("M19").Value = 25.054466230936
$ws.Range("N19").Value = -32.786885245901

# Row 20 (G.L.A.)
$ws.Range("C20").Value = 1
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 26
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = -55.932203389830
$ws.Range("M20").Value = -7.142857142857
$ws.Range("N20").Value = -91.475409836065

# Row 21 (TOTAL)
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 48
$ws.Range("E21").Value = -58.333333333333
$ws.Range("F21").Value = 128
$ws.Range("G21").Value = 172
$ws.Range("H21").Value = -25.581395348837
$ws.Range("I21").Value = 1018
$ws.Range("J21").Value = 1110
$ws.Range("K21").Value = -8.288288288288
$ws.Range("L21").Value = 27.409261576971
$ws.Range("M21").Value = 23.244552058111
$ws.Range("N21").Value = -60.280920795942

# Row 22 (Transit)
$ws.Range("C22").Value = "'0"
$ws.Range("A22").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("M22").Value = -30.769230769230

# Row 23 (Housing)
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 0
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = -9.090909090909
$ws.Range("I23").Value = 80
$ws.Range("J23").Value = 82
$ws.Range("K23").Value = -2.439024390243
$ws.Range("L23").Value = -33.333333333333
$ws.Range("M23").Value = 8.108108108108

# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 34
$ws.Range("D24").Value = 69
$ws.Range("E24").Value = -50.724637681159
$ws.Range("F24").Value = 97
$ws.Range("G24").Value = 182
$ws.Range("H24").Value = -46.703296703296
$ws.Range("I24").Value = 835
$ws.Range("J24").Value = 1312
$ws.Range("K24").Value = -36.356707317073
$ws.Range("L24").Value = 50.994575045208
$ws.Range("M24").Value = -14.358974358974

# Row 25 (Misd. Assault)
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 16.666666666666
$ws.Range("F25").Value = 37
$ws.Range("G25").Value = 47
$ws.Range("H25").Value = -21.276595744680
$ws.Range("I25").Value = 275
$ws.Range("J25").Value = 283
$ws.Range("K25").Value = -2.826855123674
$ws.Range("L25").Value = 37.5
$ws.Range("M25").Value = 2.996254681647

# Row 26 (UCR Rape*)
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = "'0"
$ws.Range("A26").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = "'***.*"
$ws.Range("A26").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("I26").Value = 23
$ws.Range("K26").Value = 4.545454545454
$ws.Range("L26").Value = 21.052631578947

# Row 27 (Other Sex Crimes)
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = -75
$ws.Range("I27").Value = 32
$ws.Range("J27").Value = 60
$ws.Range("K27").Value = -46.666666666666
$ws.Range("L27").Value = -8.571428571428

# Row 28 (Shooting Vic.)
$ws.Range("N28").Value = -85.714285714285

# Row 29 (Shooting Inc.)
$ws.Range("N29").Value = -82.352941176470

# Row 30 (Hate Crimes)
$ws.Range("F30").Value = "'0"
$ws.Range("A30").Copy()
$ws.Range("F30").PasteSpecial(-4122)

$excel.CutCopyMode = 0
Write-Host "Edit complete"
